$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.227.62"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "2.361.29"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("E4").Value = "  -0.07%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "543.49"
$r.ClearFormats()
$ws.Range("E5").Value = "  -0.12%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "133.83"
$r.ClearFormats()
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("E7").Value = "  -0.05%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.567"
$r.ClearFormats()
$ws.Range("E8").Value = "  +5.22%  "
$ws.Range("E9").Value = "  +4.45%  "
$ws.Range("E10").Value = "  +2.27%  "
$ws.Range("E11").Value = "  -2.03%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.356"
$r.ClearFormats()
$ws.Range("E12").Value = "  -0.94%  "
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "23.86"
$r.ClearFormats()
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("D14").Value = "2.778.53"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").Value = "58.146.75"
$ws.Range("E15").Value = "  +0.25%  "
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "0.0000137"
$r.ClearFormats()
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").Value = "2.359.25"
$ws.Range("E17").Value = "  +1.06%  "
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "10.85"
$r.ClearFormats()
$ws.Range("E18").Value = "  +2.44%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "4.32"
$r.ClearFormats()
$ws.Range("E19").Value = "  +2.62%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "330.60"
$r.ClearFormats()
$ws.Range("E20").Value = "  -0.86%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "6.84"
$r.ClearFormats()
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("E22").Value = "  +0.26%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "63.54"
$r.ClearFormats()
$ws.Range("E23").Value = "  +2.77%  "
$ws.Range("E24").Value = "  -1.49%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.ClearFormats()
$ws.Range("E25").Value = "  +0.11%  "
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "8.26"
$r.ClearFormats()
$ws.Range("E26").Value = "  -2.58%  "
$ws.Range("E27").Value = "  -5.99%  "
$ws.Range("E28").Value = "  -0.26%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "170.95"
$r.ClearFormats()
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("D30").Value = "0.0₃0739"
$ws.Range("E30").Value = "  +1.02%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "6.15"
$r.ClearFormats()
$ws.Range("E31").Value = "  +0.29%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "18.40"
$r.ClearFormats()
$ws.Range("E32").Value = "  -0.48%  "
$ws.Range("E33").Value = "  -2.57%  "
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "0.999"
$r.ClearFormats()
$ws.Range("E35").Value = "  -0.03%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "4.19"
$r.ClearFormats()
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("E38").Value = "  -1.86%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "0.406"
$r.ClearFormats()
$ws.Range("E39").Value = "  +7.08%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "141.93"
$r.ClearFormats()
$ws.Range("E40").Value = "  -4.77%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "3.67"
$r.ClearFormats()
$ws.Range("E41").Value = "  +1.73%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "288.80"
$r.ClearFormats()
$ws.Range("E42").Value = "  +1.57%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "0.0950"
$r.ClearFormats()
$ws.Range("E43").Value = "  +2.59%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.0519"
$r.ClearFormats()
$ws.Range("E44").Value = "  +2.68%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "18.99"
$r.ClearFormats()
$ws.Range("E45").Value = "  -1.03%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "0.567"
$r.ClearFormats()
$ws.Range("E46").Value = "  +0.90%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "0.0223"
$r.ClearFormats()
$ws.Range("E47").Value = "  +2.82%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "0.387"
$r.ClearFormats()
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("E49").Value = "  +0.18%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.951"
$r.ClearFormats()
$ws.Range("E51").Value = "  +0.59%  "
